$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("LiveData")

$ws.Range("C3").Value = 179563
$ws.Range("C4").Value = 169519
$ws.Range("C7").Value = 5.59
$ws.Range("C8").Value = 65.17
